# Scheduled runner update: refresh Kraken (market) profit figures across all Sheets
# Updates currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H:N) for the affected leve rows on each job sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 8087.5835
$ws.Cells.Item(32, 10).Value = 9999.875
$ws.Cells.Item(32, 12).Value = 9999.875
$ws.Cells.Item(32, 14).Value = -10651.875

$ws.Cells.Item(34, 8).Value = 2770.6667
$ws.Cells.Item(34, 9).Value = 1656
$ws.Cells.Item(34, 10).Value = 5000
$ws.Cells.Item(34, 11).Value = 1656
$ws.Cells.Item(34, 12).Value = 5000
$ws.Cells.Item(34, 13).Value = -1453
$ws.Cells.Item(34, 14).Value = -5406

$ws.Cells.Item(36, 8).Value = 2770.6667
$ws.Cells.Item(36, 9).Value = 1656
$ws.Cells.Item(36, 10).Value = 5000
$ws.Cells.Item(36, 11).Value = 1656
$ws.Cells.Item(36, 12).Value = 5000
$ws.Cells.Item(36, 13).Value = -941
$ws.Cells.Item(36, 14).Value = -6430

$ws.Cells.Item(62, 8).Value = 3580.8
$ws.Cells.Item(62, 9).Value = 2635
$ws.Cells.Item(62, 10).Value = 4999.5
$ws.Cells.Item(62, 11).Value = 2635
$ws.Cells.Item(62, 12).Value = 4999.5
$ws.Cells.Item(62, 13).Value = -2011
$ws.Cells.Item(62, 14).Value = -6247.5

$ws.Cells.Item(65, 8).Value = 3580.8
$ws.Cells.Item(65, 9).Value = 2635
$ws.Cells.Item(65, 10).Value = 4999.5
$ws.Cells.Item(65, 11).Value = 13175
$ws.Cells.Item(65, 12).Value = 24997.5
$ws.Cells.Item(65, 13).Value = -10055
$ws.Cells.Item(65, 14).Value = -31237.5

$ws.Cells.Item(106, 8).Value = 7801
$ws.Cells.Item(106, 9).Value = 10701.5
$ws.Cells.Item(106, 11).Value = 10701.5
$ws.Cells.Item(106, 13).Value = -10070.5

$ws.Cells.Item(132, 8).Value = 5796.875
$ws.Cells.Item(132, 9).Value = 1218.75
$ws.Cells.Item(132, 10).Value = 10375
$ws.Cells.Item(132, 11).Value = 3656.25
$ws.Cells.Item(132, 12).Value = 31125
$ws.Cells.Item(132, 13).Value = -1126.25
$ws.Cells.Item(132, 14).Value = -36185

$ws.Cells.Item(138, 8).Value = 5759.143
$ws.Cells.Item(138, 9).Value = 6634.75
$ws.Cells.Item(138, 10).Value = 4591.6665
$ws.Cells.Item(138, 11).Value = 19904.25
$ws.Cells.Item(138, 12).Value = 13774.9995
$ws.Cells.Item(138, 13).Value = -14764.25
$ws.Cells.Item(138, 14).Value = -24054.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14351.643
$ws.Cells.Item(32, 9).Value = 12094.3
$ws.Cells.Item(32, 11).Value = 12094.3
$ws.Cells.Item(32, 13).Value = -11807.3

$ws.Cells.Item(102, 8).Value = 1450
$ws.Cells.Item(102, 9).Value = 1450
$ws.Cells.Item(102, 11).Value = 1450
$ws.Cells.Item(102, 13).Value = 172

$ws.Cells.Item(109, 8).Value = 50000
$ws.Cells.Item(109, 10).Value = 50000
$ws.Cells.Item(109, 12).Value = 50000
$ws.Cells.Item(109, 14).Value = -52774

$ws.Cells.Item(122, 8).Value = 2644
$ws.Cells.Item(122, 9).Value = 2644
$ws.Cells.Item(122, 11).Value = 7932
$ws.Cells.Item(122, 13).Value = -5482

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 4821.909
$ws.Cells.Item(107, 9).Value = 5560.222
$ws.Cells.Item(107, 10).Value = 1499.5
$ws.Cells.Item(107, 11).Value = 5560.222
$ws.Cells.Item(107, 12).Value = 1499.5
$ws.Cells.Item(107, 13).Value = -3640.222
$ws.Cells.Item(107, 14).Value = -5339.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1496.5
$ws.Cells.Item(99, 9).Value = 1496.5
$ws.Cells.Item(99, 11).Value = 1496.5
$ws.Cells.Item(99, 13).Value = 1.5

$ws.Cells.Item(126, 8).Value = 1496.5
$ws.Cells.Item(126, 9).Value = 1496.5
$ws.Cells.Item(126, 11).Value = 4489.5
$ws.Cells.Item(126, 13).Value = -2019.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 200393.2
$ws.Cells.Item(4, 9).Value = 322.33334
$ws.Cells.Item(4, 11).Value = 967.0000200000001
$ws.Cells.Item(4, 13).Value = -855.0000200000001

$ws.Cells.Item(138, 8).Value = 4925
$ws.Cells.Item(138, 9).Value = 4566.6665
$ws.Cells.Item(138, 11).Value = 13699.9995
$ws.Cells.Item(138, 13).Value = -8559.999500000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 9).Value = 3751750.5
$ws.Cells.Item(14, 10).Value = 2003080
$ws.Cells.Item(14, 11).Value = 3751750.5
$ws.Cells.Item(14, 12).Value = 2003080
$ws.Cells.Item(14, 13).Value = -3751582.5
$ws.Cells.Item(14, 14).Value = -2003416

$ws.Cells.Item(80, 8).Value = 2892.6
$ws.Cells.Item(80, 9).Value = 2687.6667
$ws.Cells.Item(80, 11).Value = 2687.6667
$ws.Cells.Item(80, 13).Value = -1689.6667

$ws.Cells.Item(83, 8).Value = 2892.6
$ws.Cells.Item(83, 9).Value = 2687.6667
$ws.Cells.Item(83, 11).Value = 13438.3335
$ws.Cells.Item(83, 13).Value = -8446.333500000001

$ws.Cells.Item(102, 8).Value = 3463.7693
$ws.Cells.Item(102, 9).Value = 3548.0908
$ws.Cells.Item(102, 10).Value = 3000
$ws.Cells.Item(102, 11).Value = 3548.0908
$ws.Cells.Item(102, 12).Value = 3000
$ws.Cells.Item(102, 13).Value = -1926.0908
$ws.Cells.Item(102, 14).Value = -6244

$ws.Cells.Item(113, 8).Value = 900
$ws.Cells.Item(113, 9).Value = 900
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 900
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 1270
$ws.Cells.Item(113, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 4240.143
$ws.Cells.Item(126, 9).Value = 4970.25
$ws.Cells.Item(126, 10).Value = 3266.6667
$ws.Cells.Item(126, 11).Value = 14910.75
$ws.Cells.Item(126, 12).Value = 9800.000100000001
$ws.Cells.Item(126, 13).Value = -12440.75
$ws.Cells.Item(126, 14).Value = -14740.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 14).ClearContents()

$ws.Cells.Item(22, 8).Value = 810
$ws.Cells.Item(22, 9).Value = 734
$ws.Cells.Item(22, 10).Value = 1000
$ws.Cells.Item(22, 11).Value = 734
$ws.Cells.Item(22, 12).Value = 1000
$ws.Cells.Item(22, 13).Value = -439
$ws.Cells.Item(22, 14).Value = -1590

$ws.Cells.Item(27, 8).Value = 810
$ws.Cells.Item(27, 9).Value = 734
$ws.Cells.Item(27, 10).Value = 1000
$ws.Cells.Item(27, 11).Value = 734
$ws.Cells.Item(27, 12).Value = 1000
$ws.Cells.Item(27, 13).Value = -627
$ws.Cells.Item(27, 14).Value = -1214

$ws.Cells.Item(55, 8).Value = 2664.2222
$ws.Cells.Item(55, 9).Value = 2282.5715
$ws.Cells.Item(55, 11).Value = 2282.5715
$ws.Cells.Item(55, 13).Value = -2109.5715

$ws.Cells.Item(68, 8).Value = 2872.6155
$ws.Cells.Item(68, 9).Value = 2955.625
$ws.Cells.Item(68, 11).Value = 2955.625
$ws.Cells.Item(68, 13).Value = -2206.625

$ws.Cells.Item(71, 8).Value = 2872.6155
$ws.Cells.Item(71, 9).Value = 2955.625
$ws.Cells.Item(71, 11).Value = 14778.125
$ws.Cells.Item(71, 13).Value = -11034.125

$ws.Cells.Item(93, 8).Value = 1663.8334
$ws.Cells.Item(93, 9).Value = 1596.6
$ws.Cells.Item(93, 11).Value = 1596.6
$ws.Cells.Item(93, 13).Value = -348.5999999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(29, 8).Value = 10100.25
$ws.Cells.Item(29, 9).Value = 8750
$ws.Cells.Item(29, 11).Value = 8750
$ws.Cells.Item(29, 13).Value = -8460

$ws.Cells.Item(122, 8).Value = 225443.11
$ws.Cells.Item(122, 9).Value = 288571.16
$ws.Cells.Item(122, 11).Value = 865713.48
$ws.Cells.Item(122, 13).Value = -863263.48

$ws.Cells.Item(132, 8).Value = 9086.6875
$ws.Cells.Item(132, 10).Value = 11999.889
$ws.Cells.Item(132, 12).Value = 35999.667
$ws.Cells.Item(132, 14).Value = -41059.667
